# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45182 (2023-09-13) to 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 to row 79 in this sheet.
$ws.Range("C2:C79").Value = 45184
